$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = "resource not allocated:"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = "resources not allocated (pending for hours)"
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "NA"
$ws.Range("O4").Value = "NA"
